# Insert a new data row at row 34 (pushing the existing rows 34..62 down
# to 35..63) and populate it with the new weekly price-report entry.
# This mirrors the source diff, which inserts one new "Poroto granado"
# record ahead of the previously-first entry and leaves every other row's
# contents untouched (only their row numbers shift by +1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 34:62 down to 35:63, leaving a blank row 34 in place with
# formatting inherited from the row above (matches Excel's native
# "Insert Sheet Rows" behaviour).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record.
$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44587
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 100112030
$ws.Cells.Item(34, 7).Value = "Poroto granado"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 20
$ws.Cells.Item(34, 11).Value = 28000
$ws.Cells.Item(34, 12).Value = 28000
$ws.Cells.Item(34, 13).Value = 28000
$ws.Cells.Item(34, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(34, 16).Value = 1120
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
